$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update last_login (column J) for a few existing users ---
$ws.Cells.Item(2, 10).Value = "2026-01-07 11:37:41"

$ws.Cells.Item(4, 2).Value = "6d5a9e0fd9eb6b313c200c18f8af2904497cd4cecb9ae719c2fad9e2f71295b7"
$ws.Cells.Item(4, 10).Value = "2026-01-01 18:13:47"

$ws.Cells.Item(12, 10).Value = "2026-01-07 14:44:46"

# --- Add new user row 18: Vardaan Aggarwal ---
$ws.Cells.Item(18, 1).Value = "Vardaan"
$ws.Cells.Item(18, 3).Value = "Vardaan Aggarwal"
$ws.Cells.Item(18, 4).Value = "vardaan.aggarwal@koenig-solutions.com"
$ws.Cells.Item(18, 5).Value = "Authorised Singator and Manager in Dubai"
$ws.Cells.Item(18, 6).Value = "Admin"
$ws.Cells.Item(18, 7).Value = "EMP1636"
$ws.Cells.Item(18, 8).Value = $true
